$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old _GoBack bookmark up front (it used to sit in the Data
#    Source paragraph right before "I intend to use the data..."). We do
#    this before touching any text so it is unambiguous regardless of how
#    the engine reflows runs during Find/Replace.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
    Write-Host "Step1: removed old _GoBack bookmark"
}

# ---------------------------------------------------------------------------
# 2) Data Source paragraph: drop the "Get Featured Projects" method mention
# ---------------------------------------------------------------------------
$old1 = ". I intend to use the data obtained using the “Get All Projects” and “Get Featured Projects” methods to gather information for all of the "
$new1 = ". I intend to use the data obtained using the “Get All Projects” methods to gather information for all of the "
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Host "Step2 (Get All Projects):" $found1

# ---------------------------------------------------------------------------
# 3) Data Source paragraph: drop the "featured on website main page" clause
# ---------------------------------------------------------------------------
$old2 = " site and incorporate information on whether that project has ever been featured on website main page."
$new2 = " site."
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Host "Step3 (site.):" $found2

# ---------------------------------------------------------------------------
# 4) Techniques paragraph: insert the new visualization/statistics sentences
# ---------------------------------------------------------------------------
$old3 = "useful features. NLP techniques such as "
$new3 = "useful features. Visualization techniques will include creating maps of where projects are being run and to identify spatial patterns in projects that are funded. I’ll gather summary statistics to look for differences in funded versus retired projects.  NLP techniques such as "
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Host "Step4 (Visualization techniques):" $found3

# ---------------------------------------------------------------------------
# 5) Techniques paragraph: drop the survival-analysis sentence, leaving just
#    a trailing space in that run
# ---------------------------------------------------------------------------
$old4 = " Finally, time allowing, I’m interested in using survival analysis techniques to predict time to funding."
$new4 = " "
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
Write-Host "Step5 (drop survival analysis sentence):" $found4

# ---------------------------------------------------------------------------
# 6) Add a new _GoBack bookmark into the (now empty) paragraph that directly
#    follows the Techniques paragraph
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Techniques:")) {
        $target = $d.Paragraphs.Item($i + 1).Range
        $d.Bookmarks.Add("_GoBack", $target)
        Write-Host "Step6: added new _GoBack bookmark after paragraph" $i
        break
    }
}

# ---------------------------------------------------------------------------
# 7) Challenges paragraph: append the closing "However, ..." sentence
# ---------------------------------------------------------------------------
$old5 = "but unavailable. "
$new5 = "but unavailable. However, there is still a substantial amount of data available through the GlobalGiving API with a lot of potential to create predictive features."
$found5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
Write-Host "Step7 (However sentence):" $found5

# ---------------------------------------------------------------------------
# 8) Remove the now-superfluous empty paragraph right after the Challenges
#    paragraph (the commit collapses the blank line that used to separate
#    it from the next section)
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Challenges:")) {
        $blank = $d.Paragraphs.Item($i + 1)
        $blank.Range.Delete()
        Write-Host "Step8: removed blank paragraph after Challenges (was paragraph" ($i + 1) ")"
        break
    }
}

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
